$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.74"
$ws.Range("E2").Value = "'-3.93%"
$ws.Range("E3").Value = "'-3.90%"
$ws.Range("D4").Value = "'4.856"
$ws.Range("E4").Value = "'-2.92%"
$ws.Range("D5").Value = "'0.07167"
$ws.Range("E5").Value = "'-9.40%"
$ws.Range("D6").Value = "'7.667"
$ws.Range("E6").Value = "'-2.59%"
$ws.Range("D7").Value = "'1.731"
$ws.Range("E7").Value = "'-17.96%"
$ws.Range("D8").Value = "'3.758"
$ws.Range("E8").Value = "'-1.26%"
$ws.Range("D9").Value = "'0.8961"
$ws.Range("E9").Value = "'-3.24%"
$ws.Range("D10").Value = "'0.1648"
$ws.Range("E10").Value = "'-5.96%"
$ws.Range("D11").Value = "'0.07389"
$ws.Range("E11").Value = "'-6.61%"
$ws.Range("D12").Value = "'0.07988"
$ws.Range("E12").Value = "'-8.92%"
$ws.Range("D13").Value = "'0.03034"
$ws.Range("E13").Value = "'-4.24%"
$ws.Range("E14").Value = "'-0.67%"
$ws.Range("D15").Value = "'0.001496"
$ws.Range("E15").Value = "'-1.45%"
$ws.Range("D16").Value = "'0.005663"
$ws.Range("E16").Value = "'-5.98%"
$ws.Range("D17").Value = "'3.461"
$ws.Range("E17").Value = "'-0.18%"
$ws.Range("E18").Value = "'-7.35%"
$ws.Range("D19").Value = "'0.3291"
$ws.Range("E19").Value = "'0.12%"
$ws.Range("D20").Value = "'0.1303"
$ws.Range("E20").Value = "'0.88%"
$ws.Range("D21").Value = "'4.389"
$ws.Range("E21").Value = "'5.36%"
$ws.Range("D22").Value = "'0.2007"
$ws.Range("E22").Value = "'11.98%"
$ws.Range("D23").Value = "'0.04494"
$ws.Range("E23").Value = "'-2.57%"
$ws.Range("E24").Value = "'-1.60%"
$ws.Range("D25").Value = "'0.004007"
$ws.Range("E25").Value = "'-10.62%"
$ws.Range("E26").Value = "'0.41%"
$ws.Range("D39").Value = "'0.01636"
$ws.Range("E39").Value = "'-5.70%"
$ws.Range("D40").Value = "'0.04332"
$ws.Range("E40").Value = "'-10.14%"
$ws.Range("D41").Value = "'0.007408"
$ws.Range("E41").Value = "'0.78%"
$ws.Range("E42").Value = "'-3.66%"
$ws.Range("D43").Value = "'0.002047"
$ws.Range("E43").Value = "'-13.35%"
$ws.Range("D44").Value = "'0.01115"
$ws.Range("E44").Value = "'0.13%"
$ws.Range("D45").Value = "'0.00005739"
$ws.Range("E45").Value = "'-4.77%"
$ws.Range("D46").Value = "'0.00000000753"
$ws.Range("E46").Value = "'0.25%"
$ws.Range("D47").Value = "'2.182"
$ws.Range("E47").Value = "'165.95%"
$ws.Range("D48").Value = "'0.003011"
$ws.Range("E48").Value = "'-11.25%"
$ws.Range("D49").Value = "'0.00002108"
$ws.Range("E49").Value = "'0.25%"
$ws.Range("D50").Value = "'0.0002007"
$ws.Range("E50").Value = "'0.25%"
